$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.18"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'0.79%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'26.24"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'4.25%"
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'0.88%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.05601"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-0.22%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'6.480"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-1.23%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.8132"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.06%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.8441"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'0.01%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.06985"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'0.53%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.02840"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'0.35%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.09382"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-0.28%"
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'-0.12%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.0005958"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-93.85%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.006125"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-2.39%"
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'3.12%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.020"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.01%"
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'-1.73%"
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'-2.16%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.1333"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-0.21%"
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'0.03204"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'-1.50%"
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'-1.35%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'3.742"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-0.02%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04667"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.30%"
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'0.27%"
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'0.15%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.004572"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'0.92%"
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.00009596"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'-1.07%"
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'0.0001937"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'-0.11%"
$ws.Range("E28").ClearFormats()
$ws.Range("D40").Value = "'0.03659"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'-0.10%"
$ws.Range("E40").ClearFormats()
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1359"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-0.38%"
$ws.Range("E41").ClearFormats()
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002659"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-2.57%"
$ws.Range("E42").ClearFormats()
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003427"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-44.89%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.008258"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'0.82%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005399"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'1.90%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'-35.85%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.002429"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'20.53%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").ClearFormats()
